$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.807.83"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.750.53"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.76"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5087"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +3.78%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "40.87"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.75%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2700"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +7.27%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06211"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +4.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.754.86"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.60%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.06930"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +2.25%  "
$ws.Range("B13").Value = "Solana"
$ws.Range("C13").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.59"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +6.81%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6274"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +9.60%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.489"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.68%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "78.05"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "25.822.01"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.69"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000006725"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.93%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.979.13"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.069"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.99%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.260"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +4.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.175"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.71%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "136.69"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.20"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +4.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.457"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.89%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.768"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.35%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "102.74"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08254"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.75%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.723"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.18%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.431"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04440"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.0000"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.647"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.90%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.002"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.6060"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.60%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.692"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.961"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.51%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01563"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.66%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "101.57"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3869"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +3.42%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7545"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.923"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.99%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05508"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +7.72%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1104"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.998"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +1.55%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "30.20"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "52.91"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.70%  "
